# Reviewer-comment slide (3rd slide in the deck) lists the reviewers'
# remarks in English (shape 2) with a parallel Korean translation
# (shape 3). Items 1, 2, 3, 5 and 6 are already highlighted in
# accent6; item 4 ("Deep learning is used but the models are not
# properly explained with diagrams.") was missing that highlight in
# both the English and the Korean text. Add it there, matching the
# formatting already used by the other bullet points.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# English reviewer-comment placeholder - paragraph 5 is item "4."
$shEnglish = $s.Shapes.Item(2)
$paraEnglish = $shEnglish.TextFrame.TextRange.Paragraphs(5)
$paraEnglish.Font.Color.ObjectThemeColor = 10   # msoThemeColorAccent6

# Korean translation placeholder - paragraph 5 is item "4." too
$shKorean = $s.Shapes.Item(3)
$paraKorean = $shKorean.TextFrame.TextRange.Paragraphs(5)
$paraKorean.Font.Color.ObjectThemeColor = 10    # msoThemeColorAccent6
